$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 397, shifting existing rows 397-431 down to 398-432.
$ws.Rows.Item(397).Insert()

# Populate the newly inserted row 397 with the new record's data.
$ws.Cells.Item(397, 1).Value = 4
$ws.Cells.Item(397, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(397, 3).Value = "Los Lagos"
$ws.Cells.Item(397, 4).Value = 44578
$ws.Cells.Item(397, 5).Value = 10
$ws.Cells.Item(397, 6).Value = 100112033
$ws.Cells.Item(397, 7).Value = "Lechuga"
$ws.Cells.Item(397, 8).Value = "Escarola"
$ws.Cells.Item(397, 9).Value = "Primera"
$ws.Cells.Item(397, 10).Value = 200
$ws.Cells.Item(397, 11).Value = 13000
$ws.Cells.Item(397, 12).Value = 13000
$ws.Cells.Item(397, 13).Value = 13000
$ws.Cells.Item(397, 14).Value = "$/caja 15 unidades"
$ws.Cells.Item(397, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(397, 16).Value = 867
$ws.Cells.Item(397, 17).Value = 15
$ws.Cells.Item(397, 18).Value = "Hortaliza"
